$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.019.60"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.517.30"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.27"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.44"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.522.89"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.965.81"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.12"
$ws.Range("E15").Value = "  -1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.027.23"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.524.81"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.75"
$ws.Range("E21").Value = "  -0.50%  "

$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  +1.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.28"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("E25").Value = "  -3.05%  "

$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.79"
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.98"
$ws.Range("E32").Value = "  +2.97%  "

$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.45"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("E37").Value = "  -2.60%  "

$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.93"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("E42").Value = "  -8.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.93"
$ws.Range("E43").Value = "  -5.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.55"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.45"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0511"
$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("E51").Value = "  -1.78%  "
